# Updated cryptos list on Wed Jul 19 05:34:29 UTC 2023 with GitHub Actions
#
# Applies the latest crypto-market snapshot to the worksheet: updates prices
# and 1h volume percentages for all existing coins, and inserts a new
# "LidoStakedEther" entry (which pushes the remaining rows down by one and
# drops the last entry, "Elrond", off the bottom of the list).
#
# Each row tuple is: (row, Coin, Link, Price, Volume(1h), forceTextPrice)
# forceTextPrice is 1 when the Price string would otherwise be
# auto-interpreted as a number by Excel (e.g. "1.000" or "93.62"), so the
# cell is explicitly formatted as Text first to preserve the literal string
# exactly as published (matching the rest of the sheet, which stores every
# price as text, including values such as "30.114.54").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.097.26", "  -0.01%  ", 0),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.913.48", "  +0.33%  ", 0),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.000", "  -0.06%  ", 1),
    @(5, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.7990", "  +7.57%  ", 1),
    @(6, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "243.82", "  -0.03%  ", 1),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9999", "  -0.05%  ", 1),
    @(8, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "1.912.18", "  +0.23%  ", 0),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3195", "  +3.51%  ", 1),
    @(10, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "26.43", "  -0.23%  ", 1),
    @(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06961", "  -0.29%  ", 1),
    @(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.08000", "  -0.96%  ", 1),
    @(13, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.7534", "  -1.85%  ", 1),
    @(14, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.912.37", "  -2.25%  ", 0),
    @(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.244", "  -1.36%  ", 1),
    @(16, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "93.62", "  +1.48%  ", 1),
    @(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.113.17", "  +0.02%  ", 0),
    @(18, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "14.09", "  -1.07%  ", 1),
    @(19, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.998", "  -1.29%  ", 1),
    @(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "249.80", "  +4.10%  ", 1),
    @(21, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007831", "  +0.10%  ", 1),
    @(22, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9998", "  -0.06%  ", 1),
    @(23, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.156.46", "  -1.34%  ", 0),
    @(24, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.000", "  -0.02%  ", 1),
    @(25, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.967", "  -2.37%  ", 1),
    @(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "169.38", "  +1.41%  ", 1),
    @(27, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.336", "  -0.45%  ", 1),
    @(28, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1412", "  +11.03%  ", 1),
    @(29, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.99", "  +0.02%  ", 1),
    @(30, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.070", "  +1.10%  ", 1),
    @(31, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.388", "  +2.64%  ", 1),
    @(32, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.529", "  -0.89%  ", 1),
    @(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.361", "  +0.66%  ", 1),
    @(34, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.132", "  +1.22%  ", 1),
    @(35, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05446", "  +4.03%  ", 1),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.268", "  -2.89%  ", 1),
    @(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7406", "  -1.14%  ", 1),
    @(38, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.727", "  +0.15%  ", 1),
    @(39, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01933", "  -1.53%  ", 1),
    @(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.796", "  -0.13%  ", 1),
    @(41, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.196", "  -2.24%  ", 1),
    @(42, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4470", "  -0.60%  ", 1),
    @(43, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "73.34", "  -1.46%  ", 1),
    @(44, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.910", "  -3.35%  ", 1),
    @(45, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9999", "  -0.16%  ", 1),
    @(46, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8358", "  -0.63%  ", 1),
    @(47, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.926", "  +0.24%  ", 1),
    @(48, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.625", "  -1.45%  ", 1),
    @(49, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "100.84", "  -1.15%  ", 1),
    @(50, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "990.18", "  +6.97%  ", 1),
    @(51, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "2.065.15", "  -1.69%  ", 0)

)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    if ($row[5] -eq 1) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
}
